$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "27.300.61"
$ws.Range("E2").Value2 = "  -0.81%  "
$ws.Range("D3").Value2 = "1.705.58"
$ws.Range("E3").Value2 = "  -1.12%  "
$ws.Range("E4").Value2 = "  -0.13%  "
$ws.Range("D5").Value2 = "'223.93"
$ws.Range("E5").Value2 = "  -1.04%  "
$ws.Range("D6").Value2 = "'0.5307"
$ws.Range("E6").Value2 = "  -1.30%  "
$ws.Range("E7").Value2 = "  -0.07%  "
$ws.Range("D8").Value2 = "'0.2658"
$ws.Range("E8").Value2 = "  -1.28%  "
$ws.Range("D9").Value2 = "'0.06575"
$ws.Range("E9").Value2 = "  -0.72%  "
$ws.Range("D10").Value2 = "'20.74"
$ws.Range("E10").Value2 = "  -4.61%  "
$ws.Range("D11").Value2 = "'0.07622"
$ws.Range("E11").Value2 = "  -1.71%  "
$ws.Range("D12").Value2 = "'4.513"
$ws.Range("E12").Value2 = "  -2.95%  "
$ws.Range("B13").Value2 = "WrappedEther"
$ws.Range("C13").Value2 = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value2 = "1.724.35"
$ws.Range("E13").Value2 = "  -0.40%  "
$ws.Range("B14").Value2 = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value2 = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value2 = "1.940.23"
$ws.Range("E14").Value2 = "  -1.16%  "
$ws.Range("D15").Value2 = "'0.5769"
$ws.Range("E15").Value2 = "  -2.07%  "
$ws.Range("D16").Value2 = "0.0₅8138"
$ws.Range("E16").Value2 = "  -1.89%  "
$ws.Range("D17").Value2 = "'67.58"
$ws.Range("E17").Value2 = "  -0.88%  "
$ws.Range("D18").Value2 = "27.304.91"
$ws.Range("E18").Value2 = "  -0.90%  "
$ws.Range("D19").Value2 = "'215.24"
$ws.Range("E19").Value2 = "  -4.35%  "
$ws.Range("E20").Value2 = "  -0.06%  "
$ws.Range("E21").Value2 = "  -2.97%  "
$ws.Range("D22").Value2 = "'10.36"
$ws.Range("E22").Value2 = "  -3.48%  "
$ws.Range("D23").Value2 = "'5.953"
$ws.Range("E23").Value2 = "  -2.71%  "
$ws.Range("E24").Value2 = "  -0.17%  "
$ws.Range("D25").Value2 = "'144.17"
$ws.Range("E25").Value2 = "  -2.74%  "
$ws.Range("D26").Value2 = "'1.715"
$ws.Range("E26").Value2 = "  +1.12%  "
$ws.Range("D27").Value2 = "'0.1202"
$ws.Range("E27").Value2 = "  -2.65%  "
$ws.Range("D28").Value2 = "'7.210"
$ws.Range("E28").Value2 = "  -2.97%  "
$ws.Range("E29").Value2 = "  -4.33%  "
$ws.Range("D30").Value2 = "'0.05373"
$ws.Range("E30").Value2 = "  -3.65%  "
$ws.Range("E31").Value2 = "  -1.52%  "
$ws.Range("D32").Value2 = "'3.463"
$ws.Range("E32").Value2 = "  -3.45%  "
$ws.Range("D33").Value2 = "'3.400"
$ws.Range("E33").Value2 = "  -2.10%  "
$ws.Range("D34").Value2 = "'1.642"
$ws.Range("E34").Value2 = "  -1.74%  "
$ws.Range("D35").Value2 = "'2.867"
$ws.Range("E35").Value2 = "  +1.59%  "
$ws.Range("D36").Value2 = "'2.416"
$ws.Range("E36").Value2 = "  -1.39%  "
$ws.Range("D37").Value2 = "'0.9444"
$ws.Range("E37").Value2 = "  -2.00%  "
$ws.Range("D38").Value2 = "'0.5789"
$ws.Range("E38").Value2 = "  -2.93%  "
$ws.Range("D39").Value2 = "'0.01625"
$ws.Range("E39").Value2 = "  -1.70%  "
$ws.Range("D40").Value2 = "'5.756"
$ws.Range("E40").Value2 = "  -2.39%  "
$ws.Range("D41").Value2 = "'1.003"
$ws.Range("E41").Value2 = "  -0.10%  "
$ws.Range("B42").Value2 = "Maker"
$ws.Range("C42").Value2 = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value2 = "1.040.70"
$ws.Range("E42").Value2 = "  -1.98%  "
$ws.Range("B43").Value2 = "TrustWalletToken"
$ws.Range("C43").Value2 = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value2 = "'0.8401"
$ws.Range("E43").Value2 = "  -2.46%  "
$ws.Range("D44").Value2 = "'100.90"
$ws.Range("E44").Value2 = "  -0.83%  "
$ws.Range("D45").Value2 = "1.847.86"
$ws.Range("E45").Value2 = "  -1.10%  "
$ws.Range("E46").Value2 = "  +1.78%  "
$ws.Range("D47").Value2 = "'57.62"
$ws.Range("E47").Value2 = "  -2.53%  "
$ws.Range("D48").Value2 = "'0.4516"
$ws.Range("E48").Value2 = "  +1.89%  "
$ws.Range("E49").Value2 = "  -0.38%  "
$ws.Range("D50").Value2 = "'8.036"
$ws.Range("E50").Value2 = "  -2.52%  "
$ws.Range("E51").Value2 = "  -1.00%  "

# Reset style on cells that required a quote-prefix to stay text,
# so no visible formatting changes are introduced.
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D50").Style = "Normal"
